$d = $word.ActiveDocument

$pairs = @(
    ,@("23+14=37", "61-45=16")
    ,@("39+5=44", "29-22=7")
    ,@("88-46=42", "11-9=2")
    ,@("54-33=21", "70-36=34")
    ,@("71-44=27", "49+5=54")
    ,@("5+7=12", "96-18=78")
    ,@("79+12=91", "20+41=61")
    ,@("70-17=53", "85-66=19")
    ,@("13+44=57", "44+41=85")
    ,@("47-26=21", "77+11=88")
    ,@("80-73=7", "69+9=78")
    ,@("51-33=18", "32+35=67")
    ,@("7+80=87", "38+36=74")
    ,@("61+24=85", "85-76=9")
    ,@("40-26=14", "54+17=71")
    ,@("44-11=33", "58+1=59")
    ,@("30+61=91", "51+22=73")
    ,@("36+19=55", "30-21=9")
    ,@("27+70=97", "22+5=27")
    ,@("79-44=35", "4+19=23")
    ,@("68+0=68", "81-19=62")
    ,@("95-48=47", "97-17=80")
    ,@("72-12=60", "69-61=8")
    ,@("7+21=28", "7+30=37")
    ,@("36-34=2", "27+13=40")
    ,@("59-52=7", "80-57=23")
    ,@("48-46=2", "29-19=10")
    ,@("4+36=40", "12+80=92")
    ,@("44-43=1", "72-6=66")
    ,@("46+44=90", "97-11=86")
    ,@("26-7=19", "39+2=41")
    ,@("70-68=2", "30+62=92")
    ,@("62+14=76", "4+59=63")
    ,@("80-69=11", "89-87=2")
    ,@("72+19=91", "3+81=84")
    ,@("92-75=17", "25+33=58")
    ,@("11+38=49", "17-10=7")
    ,@("48-2=46", "98+1=99")
    ,@("38+52=90", "13+39=52")
    ,@("31+32=63", "14+10=24")
    ,@("10+2=12", "19+80=99")
    ,@("98-90=8", "44+5=49")
    ,@("21+40=61", "19+27=46")
    ,@("37+5=42", "18+67=85")
    ,@("25-0=25", "9+0=9")
    ,@("19+43=62", "82-27=55")
    ,@("89-60=29", "4+32=36")
    ,@("75+14=89", "26+8=34")
    ,@("93-9=84", "76-3=73")
    ,@("48-1=47", "91-6=85")
    ,@("52+33=85", "82-14=68")
    ,@("59+0=59", "32+50=82")
    ,@("46+49=95", "82+6=88")
    ,@("16+77=93", "39+4=43")
    ,@("86-58=28", "70-58=12")
    ,@("48-43=5", "3+6=9")
    ,@("85-8=77", "29-17=12")
    ,@("66-20=46", "55+20=75")
    ,@("37-28=9", "37-20=17")
    ,@("10+60=70", "10+39=49")
    ,@("70-60=10", "6+44=50")
    ,@("18+72=90", "28-22=6")
    ,@("10-6=4", "45-11=34")
    ,@("36-28=8", "49-38=11")
    ,@("71-23=48", "61-29=32")
    ,@("56+1=57", "9+38=47")
    ,@("43+52=95", "80+9=89")
    ,@("6+33=39", "19-14=5")
    ,@("14+24=38", "60+31=91")
    ,@("93-84=9", "87+1=88")
    ,@("11+74=85", "84-38=46")
    ,@("31+10=41", "6+70=76")
    ,@("47-22=25", "72-32=40")
    ,@("97-23=74", "18+64=82")
    ,@("22-1=21", "86-15=71")
    ,@("60-33=27", "65+7=72")
    ,@("83-56=27", "65-35=30")
    ,@("87-18=69", "74-66=8")
    ,@("50-7=43", "40-32=8")
    ,@("15-9=6", "47+14=61")
    ,@("28-26=2", "1+6=7")
    ,@("28+60=88", "16+55=71")
    ,@("69-44=25", "73-31=42")
    ,@("7+7=14", "33-4=29")
    ,@("91-29=62", "47-27=20")
    ,@("35+34=69", "4+29=33")
    ,@("75-17=58", "84-43=41")
    ,@("85+11=96", "53+7=60")
    ,@("64-21=43", "91+3=94")
    ,@("90-28=62", "82-41=41")
    ,@("23+10=33", "69+12=81")
    ,@("42+1=43", "76-57=19")
    ,@("59-44=15", "52-26=26")
    ,@("29-15=14", "92-34=58")
    ,@("89-70=19", "45-6=39")
    ,@("75-38=37", "54-34=20")
    ,@("36-10=26", "82-40=42")
    ,@("3+12=15", "27-1=26")
    ,@("78-2=76", "92+5=97")
    ,@("66-36=30", "73-35=38")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

